$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @(
    @("D2", "26.230.30"),
    @("E2", "  -0.43%  "),
    @("D3", "1.589.31"),
    @("E3", "  -0.16%  "),
    @("E4", "  -0.14%  "),
    @("D5", "211.75"),
    @("E5", "  +0.73%  "),
    @("E6", "  +0.29%  "),
    @("D7", "1.00"),
    @("D8", "0.245"),
    @("E8", "  -0.03%  "),
    @("D9", "0.0605"),
    @("E9", "  -1.09%  "),
    @("D10", "19.23"),
    @("E10", "  -2.01%  "),
    @("E11", "  +0.43%  "),
    @("D12", "1.812.52"),
    @("E12", "  -0.17%  "),
    @("D13", "1.579.21"),
    @("E13", "  -1.47%  "),
    @("D14", "4.01"),
    @("E14", "  -1.53%  "),
    @("E15", "  -0.25%  "),
    @("D16", "64.07"),
    @("E16", "  -0.94%  "),
    @("D17", "26.237.58"),
    @("E17", "  -0.43%  "),
    @("E18", "  -0.66%  "),
    @("D19", "215.16"),
    @("E19", "  +1.43%  "),
    @("E20", "  -1.94%  "),
    @("D21", "1.00"),
    @("E21", "  +0.01%  "),
    @("D22", "4.24"),
    @("E23", "  -0.62%  "),
    @("D24", "8.96"),
    @("E24", "  +0.41%  "),
    @("D25", "144.22"),
    @("E25", "  -0.49%  "),
    @("D26", "1.00"),
    @("D27", "7.00"),
    @("E27", "  -0.72%  "),
    @("E28", "  -0.52%  "),
    @("E29", "  -0.87%  "),
    @("D30", "0.0497"),
    @("E30", "  -1.66%  "),
    @("E31", "  +0.85%  "),
    @("E32", "  -0.91%  "),
    @("D33", "1.392.00"),
    @("E33", "  +7.32%  "),
    @("E34", "  -1.81%  "),
    @("E35", "  -0.30%  "),
    @("E36", "  -1.16%  "),
    @("D37", "0.584"),
    @("E37", "  -4.60%  "),
    @("E38", "  -0.66%  "),
    @("D39", "0.820"),
    @("E39", "  +0.79%  "),
    @("D40", "5.85"),
    @("E40", "  +4.14%  "),
    @("E41", "  -0.17%  "),
    @("D42", "0.769"),
    @("E42", "  +0.84%  "),
    @("E43", "  -0.14%  "),
    @("E44", "  -17.10%  "),
    @("D45", "1.724.38"),
    @("E45", "  -0.16%  "),
    @("D46", "60.96"),
    @("E46", "  -2.86%  "),
    @("D47", "86.34"),
    @("E47", "  -2.44%  "),
    @("B48", "RenderToken"),
    @("C48", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"),
    @("D48", "1.49"),
    @("E48", "  -1.96%  "),
    @("B49", "Algorand"),
    @("C49", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"),
    @("D49", "0.0982"),
    @("E49", "  -1.70%  "),
    @("B50", "Cronos"),
    @("C50", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"),
    @("D50", "0.0500"),
    @("E50", "  -1.03%  "),
    @("B51", "USDD"),
    @("C51", "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"),
    @("D51", "0.997"),
    @("E51", "  -0.29%  ")
)

foreach ($update in $cellUpdates) {
    $addr = $update[0]
    $val = $update[1]
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}
